# Applies the "new class structure" edit to data_table.xlsx:
#  - data sheet: update the two contact e-mail addresses (D2/D3), widen column D,
#    move the selection to D2 and leave the sheet no longer tab-selected
#  - mapping sheet: replace the "C" column-mapping row with a "D" row, move the
#    selection to A4, and make this sheet the active/selected tab

$wb = $excel.ActiveWorkbook

$data = $wb.Worksheets.Item("data")
$mapping = $wb.Worksheets.Item("mapping")

# --- data sheet -----------------------------------------------------------

# New e-mail addresses (shared-string text change); cell style (hyperlink look)
# is preserved automatically since only the value changes.
$data.Range("D2").Value = "lesroutes2005@yahoo.fr"
$data.Range("D3").Value = "lesroute2014@gmail.com"

# Widen column D (22.42578125 -> 30.42578125 stored "characters" width).
$data.Columns.Item(4).ColumnWidth = 29.592447916666668

# --- mapping sheet ----------------------------------------------------------

# Row 4 used to map column "C" -> location "C7"; it now maps "D" -> "C7".
$mapping.Range("A4").Value = "D"

# --- selection / active tab -------------------------------------------------

# Leave the data sheet's own selection on D2 (it stops being the tab-selected
# sheet once another sheet is selected below).
$data.Range("D2").Select()

# mapping becomes the active / tab-selected sheet, with A4 selected.
$mapping.Range("A4").Select()
